# Update gh-pages to output generated at 456a3b4
# Applies the "想去人数"(F) / "最低票价"(G) refresh across all four sheets,
# and inserts a new exhibition row (上海·KigOnly01·Kigurumi, 2024-08-11) into
# sheet "展览" ahead of the existing "第六届燃梦BACG PRO" row, shifting the
# two rows below it down by one and bumping their running index (column A).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Simple "想去人数" (F column) refreshes for the rows that are not touched
# by the insert below.
$sheet1updates = @{
    "F2"  = 1413
    "F5"  = 6396
    "F6"  = 502
    "F7"  = 1053
    "F8"  = 19
    "F9"  = 3484
    "F10" = 6713
    "F12" = 1350
    "F13" = 783
    "F14" = 103
    "F15" = 11
    "F16" = 28
    "F17" = 1124
    "F19" = 116
    "F21" = 180
    "F23" = 1008
    "F24" = 327
    "F25" = 36
    "F26" = 22
    "F27" = 114
    "F29" = 1158
    "F30" = 26
    "F31" = 74
    "F33" = 23
    "F34" = 25
    "F35" = 25
    "F36" = 332
    "F37" = 25
    "F39" = 299
    "F40" = 1177
    "F41" = 20
    "F42" = 55
    "F43" = 98
}
foreach ($key in $sheet1updates.Keys) {
    $ws1.Range($key).Value = $sheet1updates[$key]
}

# Insert a brand-new row at position 44 (pushes the old row 44 -> 45, and the
# old row 45 -> 46). Excel's default insert leaves the new row unformatted,
# so pull the row-label style (bold/centered/bordered) from the row it is
# displacing before filling in the values.
$ws1.Rows.Item(44).Insert()

$ws1.Range("A45").Copy()
$ws1.Range("A44").PasteSpecial(-4122)

$ws1.Range("A44").Value = 43
$ws1.Range("B44").Value = "'2024-08-11"
$ws1.Range("C44").Value = "上海·KigOnly01·Kigurumi"
$ws1.Range("D44").Value = "逸仙路301号靠纪念路路口 上海宝丰联大酒店"
$ws1.Range("E44").Value = "2024.08.11 10:00-08.11 17:00"
$ws1.Range("F44").Value = 0
$ws1.Range("G44").Value = 78
$ws1.Range("H44").Value = "https://show.bilibili.com/platform/detail.html?id=85291"
$ws1.Range("I44").Value = "//i1.hdslb.com/bfs/openplatform/202404/GgqdW89w1714031044427.jpeg"

# Row 45 now holds what used to be row 44 (第六届燃梦BACG PRO). Its running
# index bumps from 43 to 44, and its "想去人数" refreshes from 2 to 4.
$ws1.Range("A45").Value = 44
$ws1.Range("F45").Value = 4

# Row 46 now holds what used to be row 45 (第七届燃梦BACG PRO). Only its
# running index bumps, from 44 to 45.
$ws1.Range("A46").Value = 45

# ---------------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$sheet2updates = @{
    "F11" = 117
    "F17" = 548
    "F27" = 45
    "F30" = 707
    "F34" = 87
    "F40" = 54
}
foreach ($key in $sheet2updates.Keys) {
    $ws2.Range($key).Value = $sheet2updates[$key]
}

# ---------------------------------------------------------------------
# Sheet "本地生活" (local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3updates = @{
    "F5" = 830
    "F6" = 579
    "F8" = 1140
}
foreach ($key in $sheet3updates.Keys) {
    $ws3.Range($key).Value = $sheet3updates[$key]
}

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4updates = @{
    "F3"  = 830
    "F6"  = 579
    "F10" = 6396
    "F11" = 502
    "F12" = 1053
    "F13" = 19
    "F15" = 6713
    "F16" = 117
    "F18" = 1350
    "F23" = 548
    "F24" = 1140
    "F25" = 11
    "F29" = 116
    "F32" = 36
    "F33" = 22
    "F34" = 1158
    "F35" = 45
    "F36" = 26
    "F37" = 74
    "F41" = 25
    "F43" = 332
    "F45" = 87
    "F46" = 299
    "F49" = 98
    "F50" = 54
}
foreach ($key in $sheet4updates.Keys) {
    $ws4.Range($key).Value = $sheet4updates[$key]
}
